$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (GitHub Actions data update)
$ws.Range('D2').Value = '42.904.03'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '2.362.79'
$ws.Range('E3').Value = '  -0.97%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '315.13'
$ws.Range('E5').Value = '  -4.02%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '108.56'
$ws.Range('E6').Value = '  +8.94%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.637'
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.622'
$ws.Range('E9').Value = '  -0.60%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.33'
$ws.Range('E10').Value = '  +3.48%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0929'
$ws.Range('E11').Value = '  +0.80%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '8.62'
$ws.Range('E12').Value = '  +2.59%  '
$ws.Range('E13').Value = '  -0.83%  '
$ws.Range('E14').Value = '  +1.19%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.99'
$ws.Range('E15').Value = '  -3.19%  '
$ws.Range('D16').Value = '2.723.90'
$ws.Range('E16').Value = '  -0.39%  '
$ws.Range('D17').Value = '2.359.50'
$ws.Range('E17').Value = '  -1.05%  '
$ws.Range('D18').Value = '42.897.81'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.64'
$ws.Range('E19').Value = '  -1.61%  '
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '76.27'
$ws.Range('E21').Value = '  +1.23%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.60'
$ws.Range('E22').Value = '  -4.67%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '266.93'
$ws.Range('E23').Value = '  -1.84%  '
$ws.Range('E24').Value = '  -0.48%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.49'
$ws.Range('E25').Value = '  -6.51%  '
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.47'
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '23.32'
$ws.Range('E28').Value = '  -2.43%  '
$ws.Range('E29').Value = '  +1.82%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '36.89'
$ws.Range('E30').Value = '  +4.35%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '168.80'
$ws.Range('E31').Value = '  -2.49%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0907'
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.10'
$ws.Range('E33').Value = '  +2.92%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.91'
$ws.Range('E34').Value = '  -6.86%  '
$ws.Range('E35').Value = '  -0.58%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.118'
$ws.Range('E36').Value = '  +12.29%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.70'
$ws.Range('E37').Value = '  +1.52%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0362'
$ws.Range('E38').Value = '  +1.41%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.87'
$ws.Range('E39').Value = '  +0.27%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.68'
$ws.Range('E40').Value = '  -6.88%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '104.15'
$ws.Range('E41').Value = '  +10.62%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.51'
$ws.Range('E42').Value = '  -1.29%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.238'
$ws.Range('E43').Value = '  +4.56%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '71.26'
$ws.Range('E44').Value = '  +3.40%  '
$ws.Range('E45').Value = '  +0.14%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '12.60'
$ws.Range('E46').Value = '  +6.39%  '
$ws.Range('B47').Value = 'ordi'
$ws.Range('C47').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '81.44'
$ws.Range('E47').Value = '  +20.00%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '114.20'
$ws.Range('E48').Value = '  -1.95%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '5.57'
$ws.Range('E49').Value = '  +2.40%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.18'
$ws.Range('E50').Value = '  +1.96%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.29'
$ws.Range('E51').Value = '  +2.27%  '
